# Trade #45 closed at 2026-02-16 21:29:59 - leadlag DOWN +0.000%
#
# Applies:
#  - Closes trade #16 and #17 (rows 15/16 of the "leadlag" sheet)
#  - Appends new open trade #45 to the "leadlag" sheet (row 35)
#  - Appends the two newly closed trades (#16, #17) to "All Trades" (rows 17/18)
#  - Refreshes the "Summary" and "Comparison" aggregate sheets accordingly
#
# Notes: some text cells hold values that *look* numeric/percent/date
# ("58.8%", "3.56", "2026-02-16", ...) but must stay plain text, matching the
# source file. Setting NumberFormat = "@" (Text) before assigning .Value
# keeps Excel from silently re-interpreting them as numbers/dates/percentages.

$wb = $excel.ActiveWorkbook

function Set-TextCell($ws, $row, $col, $text) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $text
}

# ---------------------------------------------------------------------------
# Sheet: leadlag  -- close out trade #16 (row 15) and trade #17 (row 16)
# ---------------------------------------------------------------------------
$leadlag = $wb.Worksheets.Item("leadlag")

$leadlag.Cells.Item(15, 7).Value = 69720.882489            # G15 Exit Price
$leadlag.Cells.Item(15, 8).Value = "CLOSED"                 # H15 Status
$leadlag.Cells.Item(15, 9).Value = -0.6024                  # I15 P&L %
$leadlag.Cells.Item(15, 10).Value = -6.02                   # J15 P&L $
$leadlag.Cells.Item(15, 13).Value = "time_exit_5min"        # M15 Exit Reason
$leadlag.Cells.Item(15, 14).Value = 5                       # N15 Duration (min)

$leadlag.Cells.Item(16, 7).Value = 69650.458016             # G16 Exit Price
$leadlag.Cells.Item(16, 8).Value = "CLOSED"                 # H16 Status
$leadlag.Cells.Item(16, 9).Value = -0.5446                  # I16 P&L %
$leadlag.Cells.Item(16, 10).Value = -5.45                   # J16 P&L $
$leadlag.Cells.Item(16, 13).Value = "time_exit_5min"        # M16 Exit Reason
$leadlag.Cells.Item(16, 14).Value = 5                       # N16 Duration (min)

# ---------------------------------------------------------------------------
# Sheet: leadlag -- append new open trade #45 (row 35)
# ---------------------------------------------------------------------------
$leadlag.Cells.Item(35, 1).Value = 45                                  # A35 Trade #
Set-TextCell $leadlag 35 2 "2026-02-16"                                # B35 Date (looks like a date - force text)
Set-TextCell $leadlag 35 3 "21:29:59"                                  # C35 Time (looks like a time - force text)
$leadlag.Cells.Item(35, 4).Value = "leadlag"                           # D35 Strategy
$leadlag.Cells.Item(35, 5).Value = "DOWN"                              # E35 Side
$leadlag.Cells.Item(35, 6).Value = 68687.78999999999                   # F35 Entry Price
# G35 Exit Price left blank (trade is still OPEN)
$leadlag.Cells.Item(35, 8).Value = "OPEN"                              # H35 Status
$leadlag.Cells.Item(35, 9).Value = 0                                   # I35 P&L %
$leadlag.Cells.Item(35, 10).Value = 0                                  # J35 P&L $
$leadlag.Cells.Item(35, 11).Value = 0.75                               # K35 Confidence
$leadlag.Cells.Item(35, 12).Value = "Binance leading with -0.142% move" # L35 Entry Reason
# M35 Exit Reason left blank (trade is still OPEN)
$leadlag.Cells.Item(35, 14).Value = 0                                  # N35 Duration (min)

# ---------------------------------------------------------------------------
# Sheet: All Trades -- append the two newly-closed trades (#16 and #17)
# ---------------------------------------------------------------------------
$allTrades = $wb.Worksheets.Item("All Trades")

$allTrades.Cells.Item(17, 1).Value = 16                                 # A17 Trade #
Set-TextCell $allTrades 17 2 "2026-02-16"                               # B17 Date
Set-TextCell $allTrades 17 3 "21:24:38"                                 # C17 Time
$allTrades.Cells.Item(17, 4).Value = "leadlag"                          # D17 Strategy
$allTrades.Cells.Item(17, 5).Value = "DOWN"                             # E17 Side
$allTrades.Cells.Item(17, 6).Value = 69303.36500000001                  # F17 Entry Price
$allTrades.Cells.Item(17, 7).Value = 69720.882489                       # G17 Exit Price
$allTrades.Cells.Item(17, 8).Value = "CLOSED"                           # H17 Status
$allTrades.Cells.Item(17, 9).Value = -0.6024                            # I17 P&L %
$allTrades.Cells.Item(17, 10).Value = -6.02                             # J17 P&L $
$allTrades.Cells.Item(17, 11).Value = 0.75                              # K17 Confidence
$allTrades.Cells.Item(17, 12).Value = "Coinbase leading with -0.081% move" # L17 Entry Reason
$allTrades.Cells.Item(17, 13).Value = "time_exit_5min"                  # M17 Exit Reason
$allTrades.Cells.Item(17, 14).Value = 5                                 # N17 Duration (min)

$allTrades.Cells.Item(18, 1).Value = 17                                 # A18 Trade #
Set-TextCell $allTrades 18 2 "2026-02-16"                               # B18 Date
Set-TextCell $allTrades 18 3 "21:24:44"                                 # C18 Time
$allTrades.Cells.Item(18, 4).Value = "leadlag"                          # D18 Strategy
$allTrades.Cells.Item(18, 5).Value = "DOWN"                             # E18 Side
$allTrades.Cells.Item(18, 6).Value = 69273.17999999999                  # F18 Entry Price
$allTrades.Cells.Item(18, 7).Value = 69650.458016                       # G18 Exit Price
$allTrades.Cells.Item(18, 8).Value = "CLOSED"                           # H18 Status
$allTrades.Cells.Item(18, 9).Value = -0.5446                            # I18 P&L %
$allTrades.Cells.Item(18, 10).Value = -5.45                             # J18 P&L $
$allTrades.Cells.Item(18, 11).Value = 0.75                              # K18 Confidence
$allTrades.Cells.Item(18, 12).Value = "Binance leading with -0.080% move" # L18 Entry Reason
$allTrades.Cells.Item(18, 13).Value = "time_exit_5min"                  # M18 Exit Reason
$allTrades.Cells.Item(18, 14).Value = 5                                 # N18 Duration (min)

# ---------------------------------------------------------------------------
# Sheet: Summary -- refresh OVERALL and leadlag aggregate rows
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")

$summary.Cells.Item(2, 3).Value = 17                 # C2 Total Trades (OVERALL)
Set-TextCell $summary 2 4 "58.8%"                    # D2 Win Rate
Set-TextCell $summary 2 5 "+1.8833%"                 # E2 Total P&L %
Set-TextCell $summary 2 6 "+0.1108%"                 # F2 Avg Trade

$summary.Cells.Item(3, 3).Value = 33                 # C3 Total Trades (leadlag)
Set-TextCell $summary 3 4 "27.3%"                    # D3 Win Rate
Set-TextCell $summary 3 5 "+1.8380%"                 # E3 Total P&L %
Set-TextCell $summary 3 6 "+0.0557%"                 # F3 Avg Trade

# ---------------------------------------------------------------------------
# Sheet: Comparison -- refresh leadlag row
# ---------------------------------------------------------------------------
$comparison = $wb.Worksheets.Item("Comparison")

$comparison.Cells.Item(2, 2).Value = 33              # B2 Total Trades
Set-TextCell $comparison 2 3 "27.3%"                 # C2 Win Rate
Set-TextCell $comparison 2 4 "1.79"                  # D2 Profit Factor
Set-TextCell $comparison 2 6 "-0.3855%"              # F2 Avg Loss %
Set-TextCell $comparison 2 7 "1.20"                  # G2 Win/Loss Ratio
Set-TextCell $comparison 2 8 "-0.6024%"              # H2 Max Drawdown
